$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.421.58'
$ws.Range("E2").Value = '  +0.25%  '

# Row 3
$ws.Range("D3").Value = '1.620.54'
$ws.Range("E3").Value = '  +0.61%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.92'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.501'
$ws.Range("E6").Value = '  -0.09%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  -0.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("E11").Value = '  -0.74%  '

# Row 12
$ws.Range("D12").Value = '1.847.37'
$ws.Range("E12").Value = '  +0.53%  '

# Row 13
$ws.Range("D13").Value = '1.614.39'
$ws.Range("E13").Value = '  +0.48%  '

# Row 14
$ws.Range("E14").Value = '  -0.36%  '

# Row 15
$ws.Range("E15").Value = '  -0.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.82'
$ws.Range("E16").Value = '  -1.21%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '236.29'
$ws.Range("E17").Value = '  +7.13%  '

# Row 18
$ws.Range("D18").Value = '26.431.60'
$ws.Range("E18").Value = '  +0.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.85'
$ws.Range("E19").Value = '  +4.34%  '

# Row 20
$ws.Range("E20").Value = '  +0.25%  '

# Row 21
$ws.Range("E21").Value = '  +0.12%  '

# Row 22
$ws.Range("E22").Value = '  -0.98%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.20'
$ws.Range("E23").Value = '  +3.20%  '

# Row 24
$ws.Range("E24").Value = '  +0.53%  '

# Row 25
$ws.Range("E25").Value = '  +1.46%  '

# Row 26
$ws.Range("E26").Value = '  +0.09%  '

# Row 27
$ws.Range("E27").Value = '  +0.74%  '

# Row 28
$ws.Range("E28").Value = '  +0.44%  '

# Row 29
$ws.Range("E29").Value = '  +1.88%  '

# Row 30
$ws.Range("E30").Value = '  +0.04%  '

# Row 31
$ws.Range("E31").Value = '  -0.16%  '

# Row 32
$ws.Range("D32").Value = '1.518.65'
$ws.Range("E32").Value = '  +5.23%  '

# Row 33
$ws.Range("E33").Value = '  +1.44%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +0.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.51'
$ws.Range("E35").Value = '  +2.50%  '

# Row 36
$ws.Range("E36").Value = '  -0.06%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.568'
$ws.Range("E37").Value = '  +1.56%  '

# Row 38
$ws.Range("E38").Value = '  +0.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.834'
$ws.Range("E39").Value = '  -0.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.88'
$ws.Range("E40").Value = '  +0.79%  '

# Row 41
$ws.Range("E41").Value = '  +0.05%  '

# Row 42
$ws.Range("E42").Value = '  +0.88%  '

# Row 43
$ws.Range("D43").Value = '1.759.88'
$ws.Range("E43").Value = '  +0.54%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.80'
$ws.Range("E44").Value = '  +1.72%  '

# Row 45
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.762'
$ws.Range("E45").Value = '  +0.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.916'
$ws.Range("E46").Value = '  +0.73%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.60'
$ws.Range("E47").Value = '  +3.03%  '

# Row 48
$ws.Range("E48").Value = '  +1.16%  '

# Row 49
$ws.Range("E49").Value = '  -0.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0965'
$ws.Range("E50").Value = '  +0.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.49'
$ws.Range("E51").Value = '  +0.21%  '
